$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto snapshot values.
# For Price values that look like plain numbers, force text format first so Excel
# keeps the exact string (matching trailing zeros / decimal formatting) instead of
# coercing them into a floating point number.

$ws.Range("D2").Value = "37.464.83"
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("D3").Value = "2.041.25"
$ws.Range("E3").Value = "  +3.87%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.11"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.21"
$ws.Range("E7").Value = "  -1.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.393"
$ws.Range("E9").Value = "  +4.09%  "
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.16"
$ws.Range("E12").Value = "  +6.67%  "
$ws.Range("D13").Value = "2.341.46"
$ws.Range("E13").Value = "  +3.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.850"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.06"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D17").Value = "2.038.45"
$ws.Range("E17").Value = "  +3.76%  "
$ws.Range("D18").Value = "37.387.89"
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("E21").Value = "  +3.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.95"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +4.61%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.34"
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.52"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.138"
$ws.Range("E28").Value = "  -4.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.95"
$ws.Range("E29").Value = "  +3.50%  "
$ws.Range("E30").Value = "  +4.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0682"
$ws.Range("E32").Value = "  +11.06%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +11.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.51"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("E36").Value = "  +6.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.01"
$ws.Range("E40").Value = "  +3.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0982"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("E42").Value = "  +3.87%  "
$ws.Range("E43").Value = "  +1.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.74"
$ws.Range("E44").Value = "  +5.13%  "
$ws.Range("D45").Value = "1.402.95"
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.84"
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("E47").Value = "  +3.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.51"
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.11"
$ws.Range("E49").Value = "  +14.26%  "
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("D51").Value = "2.232.80"
$ws.Range("E51").Value = "  +3.84%  "
